# Weekly data refresh: insert two new rows of "Arveja Verde" price data
# (Terminal Hortofrutícola Agro Chillán) at the top of the data block,
# just below the three most-recent entries already sitting in rows 2-4.
# This pushes the existing rows 5-31 down to 7-33 and grows the used
# range from A1:R31 to A1:R33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 5 (shifts rows 5..31 down to 7..33).
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# New row 5
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44530
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112022
$ws.Range("G5").Value = "Arveja Verde"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("N5").Value = "$/saco 25 kilos"
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 580
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"

# New row 6
$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C6").Value = "Ñuble"
$ws.Range("D6").Value = 44530
$ws.Range("E6").Value = 16
$ws.Range("F6").Value = 100112022
$ws.Range("G6").Value = "Arveja Verde"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 500
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
